$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift header labels: C1, D1, E1 rotate (C<-D<-E<-C, i.e. old E1 -> C1, old C1 -> D1, old D1 -> E1)
$ws.Range("C1").Value = "bedrooms_2"
$ws.Range("D1").Value = "kitchens_2"
$ws.Range("E1").Value = "living_rooms_1"

# Keep the "selected" marker columns (C/D/E) aligned with their renamed headers
# Row 5: marker moves from D to C
$ws.Range("C5").Value = 1
$ws.Range("D5").Value = 0

# Row 6: marker moves from E to D
$ws.Range("D6").Value = 1
$ws.Range("E6").Value = 0

# Row 7: marker moves from C to E
$ws.Range("C7").Value = 0
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 1
